# Sample_Attendance_File.xlsx edit:
#  - Duplicate the "Daily Attendance Template" sheet to create a new,
#    dated sign-in sheet ("04-13-2022") placed at the front of the workbook.
#  - On that new sheet, set the "Time In" indicator cells (I2:I4) to show
#    how/when each staff member signed back in: "1:00 AM" for the night-off
#    and day-off curfews, "5:00 PM" for the leaving-camp curfew.
#  - Tweak the "# Still out of Camp" formula on the new sheet so the
#    "left camp" tally uses the "Leaving Camp" wording for staff who are
#    still out, matching the curfew label used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook

# The template sheet we clone for each day's attendance.
$template = $wb.Worksheets.Item("Daily Attendance Template")

# Copy it to before the first sheet in the workbook -> becomes the new
# first/active tab.
$template.Copy($wb.Worksheets.Item(1))

$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "04-13-2022"

# "Time In" column: indicate how staff signed back out/in for each curfew.
$newSheet.Range("I2").Value = "1:00 AM"
$newSheet.Range("I3").Value = "1:00 AM"
$newSheet.Range("I4").Value = "5:00 PM"

# Count staff still out of camp: on this sheet "Leaving Camp" is the label
# used for staff who left (vs. "Left Camp" used when they've returned).
$newSheet.Range("I8").Formula = '=(COUNTIF(E:E, "=Day Off") + COUNTIF(E:E, "=Night Off") +  COUNTIF(E:E, "=Leaving Camp")) - (COUNTIF(E1:E8, "=Day Off") + COUNTIF(E1:E8, "=Night Off") +  COUNTIF(E1:E8, "=Left Camp"))'

# Match column I to its new, shorter "h:mm AM/PM"-sized contents.
$newSheet.Columns.Item(9).AutoFit()

# Leave the selection where the user would naturally land next: the first
# empty row below the summary block.
$newSheet.Range("I9").Select()
